$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11th column); this shifts the
# former K:P data (Қарз сумма .. Реферал тел раками) one column to the
# right, into L:Q, inheriting column K's width/style.
$ws.Columns("K").Insert()

# New column header: "Тўланган сумма" (Paid amount), placed ahead of the
# "Қарз сумма" (Debt amount) column that got pushed to L.
$ws.Range("K4").Value = "Тўланган сумма"

# Update the hidden AutoFilter database range so it covers the new
# rightmost column (was A4:P4, now A4:Q4).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "TDSheet!_FilterDatabase") {
        $n.RefersTo = "=TDSheet!`$A`$4:`$Q`$4"
    }
}

# Reflect the new selection / view state (column K selected, no frozen
# top-left scroll offset).
$ws.Columns("K:K").Select()
